# Day-5 - JavaScript basics: Scope var,let, const keywords covered
# Mark the "Chapter 14: Scope - var,let, const" row (row 16) as DONE,
# matching the style already used for the other completed rows (B2:B15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "DONE" cell format (green theme fill, style index reused)
# from B2 and apply it to B16 before writing the value, so B16 ends up
# sharing the exact same cell style as the rest of the STATUS column.
$ws.Range("B2").Copy()
$ws.Range("B16").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B16").Value = "DONE"

# Reflect the scrolled viewport / new active selection from the edit session.
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
$ws.Range("D16").Select()
